$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the typo "workspaceer_description" before touching anything.
# ------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute("workspaceer_description", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$textStart = $found.Start
$textEnd   = $found.End

# ------------------------------------------------------------------
# 2. Split the run at the start of the typo (between "  def " and
#    "workspaceer_description") using a throw-away bookmark. Doing this
#    *before* the text is rewritten keeps the untouched "  def " run from
#    inheriting formatting quirks off the text that is about to change,
#    and later lets the fixed word live in its own run exactly like the
#    target markup: one run for "  def " and a second, distinct run for
#    "modeler_description".
# ------------------------------------------------------------------
$splitRange = $d.Range($textStart, $textStart)
$d.Bookmarks.Add("_TempSplit", $splitRange)

# ------------------------------------------------------------------
# 3. Fix the typo: "workspaceer_description" -> "modeler_description".
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("workspaceer_description", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "modeler_description", 2)
$newTextEnd = $rng.End

# Drop the helper bookmark now that the edit is in place.
$d.Bookmarks("_TempSplit").Delete()

# ------------------------------------------------------------------
# 4. Move "_GoBack" so it sits right after the freshly typed
#    "modeler_description" (mirroring Word dropping "_GoBack" at the
#    location of the last edit). The previous "_GoBack" - which lived on
#    "SetpointManager:OutdoorAirReset" - is implicitly replaced because
#    bookmark names must stay unique, but it is also removed explicitly
#    first to be safe.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# A zero-length range sitting exactly on a paragraph's trailing edge (i.e.
# right before the paragraph mark) can't be targeted directly and ends up
# snapping to the start of the document, so a temporary placeholder
# character is inserted right after the new text, "_GoBack" is anchored
# just in front of it, and the placeholder is deleted again.
$placeholderRange = $d.Range($newTextEnd, $newTextEnd)
$placeholderRange.InsertAfter("Z")

$goBackRange = $d.Range($newTextEnd, $newTextEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$placeholderRange2 = $d.Range($newTextEnd, $newTextEnd + 1)
$placeholderRange2.Text = ""
